$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "street_classifier_2021_12_03_12_30_57"
$ws.Range("B6").Value = "Window size 50, rbf in SVC classifier"

$ws.Range("B6").Select()
